# Edit script for czech-republic_cfl-group-a_2023-2024 workbook
# - Swaps/rotates the F:V ("match detail") content of several row-pairs/groups
#   (home/away team names, scores, odds, timestamps, urls), while leaving
#   columns A-E (Indice, pais, torneio, temporada, data_partida) untouched.
# - Appends 5 brand-new match rows (61-65) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fvCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowFV($ws, $row) {
    $vals = @()
    foreach ($c in $fvCols) {
        $vals += $ws.Range("$c$row").Value()
    }
    return $vals
}

function Set-RowFV($ws, $row, $vals) {
    for ($i = 0; $i -lt $fvCols.Length; $i++) {
        $ws.Range("$($fvCols[$i])$row").Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 1. Re-order the "match detail" columns (F:V) across certain rows.
#    Columns A:E (Indice/pais/torneio/temporada/data_partida) stay in place.
# ---------------------------------------------------------------------------

# Simple two-way swaps.
$swapPairs = @(
    @(10, 11),
    @(14, 15),
    @(16, 17),
    @(43, 44),
    @(48, 49)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $v1 = Get-RowFV $ws $r1
    $v2 = Get-RowFV $ws $r2
    Set-RowFV $ws $r1 $v2
    Set-RowFV $ws $r2 $v1
}

# Three-way rotation: row 29 <- old row31, row 30 <- old row29, row 31 <- old row30.
$v29 = Get-RowFV $ws 29
$v30 = Get-RowFV $ws 30
$v31 = Get-RowFV $ws 31
Set-RowFV $ws 29 $v31
Set-RowFV $ws 30 $v29
Set-RowFV $ws 31 $v30

# ---------------------------------------------------------------------------
# 2. Append 5 new rows (61-65) at the end of the sheet.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Indice=60; E=45192.6875;         F="Domazlice";         G=4; H="Motorlet Prague";  I=0; J=1.33; K="22/09/2023 03:43"; L=1.17; M="23/09/2023 16:26"; N=4.8;  O="22/09/2023 03:43"; P=6.83;  Q="23/09/2023 16:26"; R=5.95; S="22/09/2023 03:43"; T=11.93; U="23/09/2023 16:26"; V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/domazlice-motorlet-prague/db2Gg1tA/" },
    @{ Indice=61; E=45193.42708333334;  F="Admira Prague";     G=3; H="Plzen B";           I=2; J=2.64; K="22/09/2023 21:42"; L=2.74; M="24/09/2023 08:16"; N=3.35; O="22/09/2023 21:42"; P=3.39;  Q="24/09/2023 08:18"; R=2.22; S="22/09/2023 21:42"; T=2.36;  U="24/09/2023 08:16"; V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/admira-prague-plzen/4O5OiuBM/" },
    @{ Indice=62; E=45193.42708333334;  F="Slavia Prague B";   G=5; H="Dukla Prague B";    I=1; J=1.18; K="22/09/2023 21:42"; L=1.14; M="24/09/2023 09:25"; N=6.65; O="22/09/2023 21:42"; P=8.5;   Q="24/09/2023 09:30"; R=7.79; S="22/09/2023 21:42"; T=11.59; U="24/09/2023 09:30"; V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/slavia-prague-dukla-prague/nmeTjaQS/" },
    @{ Indice=63; E=45193.4375;         F="Taborsko akademie"; G=1; H="Bohemians 1905 B";  I=0; J=2.49; K="23/09/2023 03:42"; L=1.81; M="24/09/2023 10:13"; N=3.58; O="23/09/2023 03:42"; P=4.18;  Q="24/09/2023 10:13"; R=2.24; S="23/09/2023 03:42"; T=3.06;  U="24/09/2023 10:13"; V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/taborsko-akademie-bohemians-1905/O8A3d3di/" },
    @{ Indice=64; E=45193.64583333334;  F="Hostoun";           G=2; H="Karlovy Vary";      I=0; J=1.64; K="23/09/2023 03:42"; L=1.68; M="24/09/2023 14:49"; N=3.87; O="23/09/2023 03:42"; P=4.04;  Q="24/09/2023 14:49"; R=3.85; S="23/09/2023 03:42"; T=4.15;  U="24/09/2023 14:49"; V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/hostoun-karlovy-vary/zqE7eNBc/" }
)

$destRow = 61
foreach ($rowData in $newRows) {

    # Copy the whole-row formatting (bold/centered style on A, date style on E)
    # from the last existing data row so new cells match existing look & feel.
    $ws.Range("A60").Copy()
    $ws.Range("A$destRow").PasteSpecial(-4122)
    $ws.Range("E60").Copy()
    $ws.Range("E$destRow").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("A$destRow").Value = $rowData.Indice
    $ws.Range("B$destRow").Value = "czech-republic"
    $ws.Range("C$destRow").Value = "cfl-group-a"
    $ws.Range("D$destRow").Value = "2023-2024"
    $ws.Range("E$destRow").Value = $rowData.E
    $ws.Range("F$destRow").Value = $rowData.F
    $ws.Range("G$destRow").Value = $rowData.G
    $ws.Range("H$destRow").Value = $rowData.H
    $ws.Range("I$destRow").Value = $rowData.I
    $ws.Range("J$destRow").Value = $rowData.J
    $ws.Range("K$destRow").Value = $rowData.K
    $ws.Range("L$destRow").Value = $rowData.L
    $ws.Range("M$destRow").Value = $rowData.M
    $ws.Range("N$destRow").Value = $rowData.N
    $ws.Range("O$destRow").Value = $rowData.O
    $ws.Range("P$destRow").Value = $rowData.P
    $ws.Range("Q$destRow").Value = $rowData.Q
    $ws.Range("R$destRow").Value = $rowData.R
    $ws.Range("S$destRow").Value = $rowData.S
    $ws.Range("T$destRow").Value = $rowData.T
    $ws.Range("U$destRow").Value = $rowData.U
    $ws.Range("V$destRow").Value = $rowData.V

    $destRow++
}

Write-Output "Done."
